# updated task used in testing
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Revised data values in row 2
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 46

# Leave the selection on C2 (as last left by the editor)
$ws.Range("C2").Select()
